# Update "想去人数" (want-to-go count) figures in column F across sheets,
# matching the data refresh recorded in the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3643
$ws1.Range("F5").Value = 3643
$ws1.Range("F6").Value = 272
$ws1.Range("F7").Value = 5170
$ws1.Range("F8").Value = 548
$ws1.Range("F9").Value = 378
$ws1.Range("F10").Value = 205
$ws1.Range("F11").Value = 704
$ws1.Range("F13").Value = 103
$ws1.Range("F14").Value = 38
$ws1.Range("F22").Value = 4949
$ws1.Range("F26").Value = 6075
$ws1.Range("F29").Value = 3232
$ws1.Range("F30").Value = 349
$ws1.Range("F35").Value = 144
$ws1.Range("F36").Value = 1059
$ws1.Range("F37").Value = 84
$ws1.Range("F40").Value = 885
$ws1.Range("F41").Value = 1039
$ws1.Range("F42").Value = 2038

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1128

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1128
$ws4.Range("F7").Value = 3643
$ws4.Range("F8").Value = 3643
$ws4.Range("F9").Value = 272
$ws4.Range("F10").Value = 5170
$ws4.Range("F11").Value = 548
$ws4.Range("F12").Value = 378
$ws4.Range("F13").Value = 205
$ws4.Range("F14").Value = 704
$ws4.Range("F16").Value = 103
$ws4.Range("F17").Value = 38
$ws4.Range("F26").Value = 4949
$ws4.Range("F30").Value = 6075
$ws4.Range("F33").Value = 3232
$ws4.Range("F34").Value = 349
$ws4.Range("F40").Value = 144
$ws4.Range("F41").Value = 1059
$ws4.Range("F42").Value = 84
$ws4.Range("F45").Value = 885
$ws4.Range("F46").Value = 1039
$ws4.Range("F48").Value = 2038
